$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Practice 1" day-of-week entry (column L) is being reset to "N/A"
# (the author re-did the warmup problems, clearing the stale day markers).
$rows = @(3,4,5,6,7,9,10,11,13,14,15,16,19,21,22,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39)
foreach ($r in $rows) {
    $ws.Range("L$r").Value = "N/A"
}

# Rows 4 and 5 also had their "Practice 2" column (M) switched from the old
# day-of-week text to a numeric placeholder.
$ws.Range("M4").Value = 2
$ws.Range("M5").Value = 2

# Radix Sort (row 34) gets its Big-O Space / Big-O Time / Other concerns filled in.
# Order matters so new shared strings land at the same indices as the source edit.
$ws.Range("G34").Value = "O(nk) where k is the length of the longest number"
$ws.Range("H34").Value = "Can only use with integers"
$ws.Range("F34").Value = "O(k + n)"

# Move the active selection to the cell that was last edited.
[void]$ws.Range("F34").Select()
